# March 24 update 3
# Add three new trailing columns (renewd, PlanID, iteration) to Sheet1,
# filling header row 1 (M1:O1) and all data rows 2:33 (M:O) with the
# "after" snapshot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells (M1:O1) the same look (bold / border / etc.)
# as the existing header cells by copying the format from L1.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Populate the new data columns for every existing data row (2-33).
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"
    $ws.Cells.Item($r, 14).Value = 20141190
    $ws.Cells.Item($r, 15).Value = 10
}
